# Updated symbol list on Thu Feb 16 16:54:03 UTC 2023 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) values to the cryptos sheet,
# keeping the cells as plain text (matching the original inlineStr storage)
# rather than letting Excel auto-convert numeric-looking strings to numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    # Leading apostrophe forces Excel to treat the entry as text even
    # when it looks numeric/percentage-like; resetting the style back to
    # "Normal" afterwards avoids leaving a stray @ (Text) number format
    # on the cell, matching the source file which carries no style index
    # on these data cells.
    $Range.Value = "'" + $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '321.42'
Set-TextValue $ws.Range("E2") '6.57%'
Set-TextValue $ws.Range("D3") '49.26'
Set-TextValue $ws.Range("E3") '12.10%'
Set-TextValue $ws.Range("D4") '5.353'
Set-TextValue $ws.Range("E4") '5.51%'
Set-TextValue $ws.Range("D5") '0.08055'
Set-TextValue $ws.Range("E5") '4.60%'
Set-TextValue $ws.Range("D6") '4.605'
Set-TextValue $ws.Range("E6") '4.26%'
Set-TextValue $ws.Range("D7") '1.397'
Set-TextValue $ws.Range("E7") '34.14%'
Set-TextValue $ws.Range("D8") '1.637'
Set-TextValue $ws.Range("E8") '1.29%'
Set-TextValue $ws.Range("D9") '0.1277'
Set-TextValue $ws.Range("E9") '0.19%'
Set-TextValue $ws.Range("D10") '0.1966'
Set-TextValue $ws.Range("E10") '5.13%'
Set-TextValue $ws.Range("D11") '0.09605'
Set-TextValue $ws.Range("E11") '4.54%'
Set-TextValue $ws.Range("D12") '0.04669'
Set-TextValue $ws.Range("E12") '11.96%'
Set-TextValue $ws.Range("E13") '-0.52%'
Set-TextValue $ws.Range("D14") '0.001317'
Set-TextValue $ws.Range("E14") '2.58%'
Set-TextValue $ws.Range("D15") '0.04181'
Set-TextValue $ws.Range("E15") '-0.27%'
Set-TextValue $ws.Range("D16") '0.005783'
Set-TextValue $ws.Range("E16") '0.45%'
Set-TextValue $ws.Range("D17") '3.343'
Set-TextValue $ws.Range("E17") '-0.01%'
Set-TextValue $ws.Range("D18") '2.459'
Set-TextValue $ws.Range("E18") '5.51%'
Set-TextValue $ws.Range("D19") '0.3506'
Set-TextValue $ws.Range("E19") '4.60%'
Set-TextValue $ws.Range("D20") '8.097'
Set-TextValue $ws.Range("E20") '-6.61%'
Set-TextValue $ws.Range("D21") '0.1371'
Set-TextValue $ws.Range("E21") '-2.02%'
Set-TextValue $ws.Range("E22") '-2.71%'
Set-TextValue $ws.Range("D23") '0.001314'
Set-TextValue $ws.Range("E23") '2.23%'
Set-TextValue $ws.Range("D24") '0.004305'
Set-TextValue $ws.Range("E24") '-3.75%'
Set-TextValue $ws.Range("D25") '0.0001348'
Set-TextValue $ws.Range("E25") '-0.06%'
Set-TextValue $ws.Range("D26") '0.0003528'
Set-TextValue $ws.Range("E26") '-95.25%'
Set-TextValue $ws.Range("D38") '0.02746'
Set-TextValue $ws.Range("E38") '10.00%'
Set-TextValue $ws.Range("D39") '0.06174'
Set-TextValue $ws.Range("E39") '16.55%'
Set-TextValue $ws.Range("D40") '0.01046'
Set-TextValue $ws.Range("E40") '76.38%'
Set-TextValue $ws.Range("D41") '0.008062'
Set-TextValue $ws.Range("E41") '4.96%'
Set-TextValue $ws.Range("D42") '0.1461'
Set-TextValue $ws.Range("E42") '8.38%'
Set-TextValue $ws.Range("D43") '0.007880'
Set-TextValue $ws.Range("E43") '7.04%'
Set-TextValue $ws.Range("D44") '0.008638'
Set-TextValue $ws.Range("E44") '14.31%'
Set-TextValue $ws.Range("D45") '0.3496'
Set-TextValue $ws.Range("E45") '16.53%'
Set-TextValue $ws.Range("D46") '0.00006632'
Set-TextValue $ws.Range("E46") '-0.31%'
Set-TextValue $ws.Range("D47") '0.00000000749'
Set-TextValue $ws.Range("E47") '-0.06%'
Set-TextValue $ws.Range("D48") '0.05504'
Set-TextValue $ws.Range("E48") '32.45%'
Set-TextValue $ws.Range("D49") '0.003988'
Set-TextValue $ws.Range("E49") '-5.07%'
Set-TextValue $ws.Range("D50") '0.00002097'
Set-TextValue $ws.Range("E50") '-0.06%'
Set-TextValue $ws.Range("D51") '0.0001997'
Set-TextValue $ws.Range("E51") '-0.06%'
